# Rename the "Documentation" tab to "Help" and the "DesireEditedSequences"
# tab to "TargetedSearch" (tab renamed as part of the
# amplicount_sequences -> amplicount_config_tsearch /
# amplicount_desired_edits -> amplicount_tsearch rework).
$wb = $excel.ActiveWorkbook

$wsHelp = $wb.Worksheets.Item("Documentation")
$wsHelp.Name = "Help"

$wsTargetedSearch = $wb.Worksheets.Item("DesireEditedSequences")
$wsTargetedSearch.Name = "TargetedSearch"

# Update the saved cursor/selection on a couple of sheets and move the
# "active" tab from Amplicon over to Help.
$wsAmplicon = $wb.Worksheets.Item("Amplicon")
$wsAmplicon.Activate() | Out-Null
$wsAmplicon.Range("A2").Select() | Out-Null

$wsLayout = $wb.Worksheets.Item("Layout")
$wsLayout.Activate() | Out-Null
$wsLayout.Range("C2").Select() | Out-Null

$wsHelp.Activate() | Out-Null
$wsHelp.Range("A3:C3").Select() | Out-Null
